$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before current row 3 (Bord Gáis - Smart EV Bonus),
# shifting everything down, to add "Yuno Energy - D Smart Bonus + 6%".
$ws.Rows.Item(3).Insert()
$ws.Range("A3").Value = "Yuno Energy - D Smart Bonus + 6%"
$ws.Range("B3").Value = 1142.77

# After the first insert, "Flogas - Smart Bonus" is now row 5 and
# "Flogas - Fixed" is row 6. Insert a new row before row 6 to add
# "Energia Offer" between them.
$ws.Rows.Item(6).Insert()
$ws.Range("A6").Value = "Energia Offer"
$ws.Range("B6").Value = 1311.46
